$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Unprotect sheet to allow editing locked cells (we'll re-protect at the end)
$ws.Unprotect()

$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1

$ws.Range("K3").Value = 1

$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1

$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("M5").Value = 11

$ws.Range("K6").Value = 1
$ws.Range("M6").Value = "+"

$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1
$ws.Range("M7").Value = "+"

$ws.Range("M3").Select()

$ws.Protect()
